$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TDCNST1")

# Remove the leftover "Pass"/"Pass : 45"/"Pass : 18"/"1.00" helper values that
# were left in row 2 (columns AH:AO and AQ:AS) - these are no longer needed.
$ws.Range("AH2:AO2").ClearContents()
$ws.Range("AQ2:AS2").ClearContents()

# Update the active selection to reflect where the user ended up working (AQ10).
$ws.Range("AQ10").Select()
